{"js": "// Capitalize the word \"order\" -> \"Order\" at the start of the two\n// \"Order\" cross-reference sentences in the MINTI order template:\n//   \"Where this order, or any rule or practice direction, ...\"\n//   \"This order has been made without hearing. ...\"\n// (The already-correct \"this Order set aside or varied...\" sentence,\n// and the all-caps \"ORDER\"/\"IT IS ORDERED THAT:\" headings, must stay\n// untouched.)\n\nconst body = context.document.body;\n\n// Case-sensitive search for the lower-case word \"order\" finds exactly\n// the two occurrences that need to become \"Order\"; it does not match\n// the already-capitalised \"Order\" later in the document.\nconst results = body.search(\"order\", { matchCase: true, matchWholeWord: true });\nresults.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  const found = results.items[i];\n  if (found.text === \"order\") {\n    // Replacing in place preserves each run's existing formatting\n    // (e.g. the bold run in \"This order has been made...\").\n    found.insertText(\"Order\", Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Capitalize the word \"order\" -> \"Order\" at the start of the two\n# \"Order\" cross-reference sentences in the MINTI order template:\n#   \"Where this order, or any rule or practice direction, ...\"\n#   \"This order has been made without hearing. ...\"\n# (The already-correct \"this Order set aside or varied...\" sentence,\n# and the all-caps \"ORDER\"/\"IT IS ORDERED THAT:\" headings, must stay\n# untouched.)\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"order\"\n$find.Replacement.Text = \"Order\"\n$find.Forward = $true\n$find.Wrap = 0\n\n# Case-sensitive, whole-word match so only the lower-case \"order\"\n# occurrences are hit (leaves \"ORDER\", \"ORDERED\" and the already\n# capitalised \"Order\" untouched), replacing every remaining run's\n# formatting preserved as-is.\n$find.Execute(\n    [ref]\"order\",\n    [ref]$true,\n    [ref]$true,\n    [ref]$false,\n    [ref]$false,\n    [ref]$false,\n    [ref]$true,\n    [ref]0,\n    [ref]$false,\n    [ref]\"Order\",\n    [ref]2\n)\n"}
